$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('C2').Value = '555 Robson St, VANCOUVER, BC, V6B 1A6'

$ws.Range('C3').Value = '4601 Westway Park Boulevard, Houston, TX, 77041'

$ws.Range('C4').Value = '#3 CALLE ACACIA OFICINA 201-B, SAN JUAN, PR, 920'
$ws.Range('D4').Value = 18.46633
$ws.Range('E4').Value = -66.10572

$ws.Range('C5').Value = '171 AUHANA RD, KIHEI, HI, 96753'

$ws.Range('C6').Value = '98-723 KUAHAO PL STE A13, PEARL CITY, HI, 96782'

$ws.Range('C7').Value = '16-711 MILO ST UNIT B, KEAAU, HI, 96749'
$ws.Range('D7').Value = 19.627332
$ws.Range('E7').Value = -155.030647

$ws.Range('C8').Value = '3989 1ST ST STE E, LIVERMORE, CA, 94551'
$ws.Range('D8').Value = 37.6781
$ws.Range('E8').Value = -121.784928

$ws.Range('C9').Value = '10604 N PALISADES WAY, BOISE, ID, 83714'

$ws.Range('C10').Value = '6344 US RT 22-3, MORROW, OH, 45152'
$ws.Range('D10').Value = 40.509237
$ws.Range('E10').Value = -82.794097

$ws.Range('C11').Value = '1725 N PACKER RD, SPRINGFIELD, MO, 65803'

$ws.Range('C12').Value = '932 KENMORE AVE, BUFFALO, NY, 14216-1462'

$ws.Range('C13').Value = '5510 N. HWY 27, BRYANT, IN, 47326-8835'
$ws.Range('D13').Value = 35.74174
$ws.Range('E13').Value = -91.65208

$ws.Range('C14').Value = '7511 KEKAA ST, HONOLULU, HI, 96825-2805'
$ws.Range('D14').Value = 21.40572
$ws.Range('E14').Value = -157.789396

$ws.Range('C15').Value = '3260 OLD FARM LN, COMMERCE, MI, 48390'

$ws.Range('C16').Value = '6150 W CHANDLER BLVD #17, CHANDLER, AZ, 85226'
$ws.Range('D16').Value = 33.310045
$ws.Range('E16').Value = -111.861363

$ws.Range('C17').Value = '680 Redna Terrace, Cincinnati, OH, 45215'

$ws.Range('C18').Value = '91-6221 KAPOLEI PARKWAY UNIT 11, EWA BEACH, HI, 96706'
$ws.Range('D18').Value = 21.325072
$ws.Range('E18').Value = -158.028212

$ws.Range('C19').Value = '35 OWOSSO DR, EUGENE, OR, 97404-2628'

$ws.Range('C20').Value = '14500 ROSCOE BLVD 4TH FLOOR, PANORAMA CITY, CA, 91402'
$ws.Range('D20').Value = 34.221238
$ws.Range('E20').Value = -118.444706

$ws.Range('C21').Value = '127 W MAIN ST APT A, JACKSON, MO, 63755-1879'
$ws.Range('D21').Value = 37.383607
$ws.Range('E21').Value = -89.677479

$ws.Range('C22').Value = '132 TERRACE DR, INDEPENDENCE, IA, 50644'
$ws.Range('D22').Value = 42.467132
$ws.Range('E22').Value = -91.878241

$ws.Range('C23').Value = '535 PINE ST, CENTRAL FALLS, RI, 2863'
$ws.Range('D23').Value = 41.884462
$ws.Range('E23').Value = -71.396381

$ws.Range('C24').Value = '6210 CARDWELL RD, CORRYTON, TN, 37721-3715'
$ws.Range('D24').Value = 36.15369
$ws.Range('E24').Value = -83.78241

$ws.Range('C25').Value = '801 PRESSLEY RD STE 100-C, CHARLOTTE, NC, 28217'
$ws.Range('D25').Value = 35.188483
$ws.Range('E25').Value = -80.893687

$ws.Range('C26').Value = '12914 SE 257TH ST, KENT, WA, 98030'

$ws.Range('C27').Value = '4363 BURTON LN, NORTH GARDEN, VA, 22959'
$ws.Range('D27').Value = 37.9407
$ws.Range('E27').Value = -78.63668

$ws.Range('C28').Value = '10415 HARMON RD, BERLIN HEIGHTS, OH, 44814'
$ws.Range('D28').Value = 41.323565
$ws.Range('E28').Value = -82.490209

$ws.Range('C29').Value = '841 WATSON LN W, NEW BRAUNFELS, TX, 78130'
$ws.Range('D29').Value = 29.694703
$ws.Range('E29').Value = -98.116089

$ws.Range('C30').Value = '903 S DEER RUN, ELLETTSVILLE, IN, 47429'

$ws.Range('C31').Value = '7735 WINTON DR, INDIANAPOLIS, IN, 46268'

$ws.Range('C32').Value = '9285 UPSTREAM LN, KNOXVILLE, TN, 37931'
$ws.Range('D32').Value = 35.973011
$ws.Range('E32').Value = -83.969472

$ws.Range('C33').Value = '930 BONNIE LN, AUBURN, CA, 95603-9452'
$ws.Range('D33').Value = 38.898076
$ws.Range('E33').Value = -121.071247

$ws.Range('C34').Value = '138 STERLING AVE, RITTMAN, OH, 44270-1655'
$ws.Range('D34').Value = 40.966545
$ws.Range('E34').Value = -81.788777

$ws.Range('C35').Value = '180 GRANTON EDGE LN, SUMMERVILLE, SC, 29486'
$ws.Range('D35').Value = 33.012176
$ws.Range('E35').Value = -80.182766

$ws.Range('C36').Value = '9756 TALL TIMBER DR, WEST CHESTER, OH, 45241-1221'
$ws.Range('D36').Value = 39.306282
$ws.Range('E36').Value = -84.381606

$ws.Range('C37').Value = '5908 TRIANGLE DR, RALEIGH, NC, 27617'

$ws.Range('C38').Value = '1132 E MARKET ST BAY 5, CHARLOTTESVILLE, VA, 22902-5351'
$ws.Range('D38').Value = 38.036141
$ws.Range('E38').Value = -78.482369

$ws.Range('C39').Value = '926 E CHURCH ST, SANDWICH, IL, 60548'
$ws.Range('D39').Value = 41.644714
$ws.Range('E39').Value = -88.616303
